# Auto-generated edit script applying the Phantom_Profits market-data refresh diff.
# Each hunk corresponds to one leve row; columns H-N hold raw market-board snapshot
# numbers (currentAveragePrice*, LevePrice*, LeveProfit*) with no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2935.8333
$ws.Range("I80").Value = 1432.5
$ws.Range("J80").Value = 5942.5
$ws.Range("K80").Value = 4297.5
$ws.Range("L80").Value = 17827.5
$ws.Range("M80").Value = -3299.5
$ws.Range("N80").Value = -19823.5
$ws.Range("H83").Value = 2935.8333
$ws.Range("I83").Value = 1432.5
$ws.Range("J83").Value = 5942.5
$ws.Range("K83").Value = 12892.5
$ws.Range("L83").Value = 53482.5
$ws.Range("M83").Value = -7900.5
$ws.Range("N83").Value = -63466.5
$ws.Range("H86").Value = 2428.4707
$ws.Range("J86").Value = 1187.8334
$ws.Range("L86").Value = 1187.8334
$ws.Range("N86").Value = -3433.8334
$ws.Range("H89").Value = 2428.4707
$ws.Range("J89").Value = 1187.8334
$ws.Range("L89").Value = 5939.166999999999
$ws.Range("N89").Value = -17171.167
$ws.Range("H138").Value = 8299.387000000001
$ws.Range("I138").Value = 6061.4116
$ws.Range("J138").Value = 11016.929
$ws.Range("K138").Value = 18184.2348
$ws.Range("L138").Value = 33050.787
$ws.Range("M138").Value = -13044.2348
$ws.Range("N138").Value = -43330.787

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1420
$ws.Range("I45").Value = 1525
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1525
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -1148
$ws.Range("N45").Value = -1754
$ws.Range("H63").Value = 1100
$ws.Range("I63").Value = 1100
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1100
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("N63").Value = -414
$ws.Range("H66").Value = 1100
$ws.Range("I66").Value = 1100
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 5500
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = -2068
$ws.Range("H96").Value = 16675.4
$ws.Range("J96").Value = 16675.4
$ws.Range("L96").Value = 16675.4
$ws.Range("N96").Value = -22167.4
$ws.Range("H130").Value = 73114.39999999999
$ws.Range("J130").Value = 73114.39999999999
$ws.Range("L130").Value = 73114.39999999999
$ws.Range("N130").Value = -83154.39999999999
$ws.Range("H131").Value = 83000
$ws.Range("J131").Value = 83000
$ws.Range("L131").Value = 83000
$ws.Range("N131").Value = -93080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 37060284
$ws.Range("I86").Value = 34033
$ws.Range("K86").Value = 34033
$ws.Range("M86").Value = -32910
$ws.Range("H89").Value = 37060284
$ws.Range("I89").Value = 34033
$ws.Range("K89").Value = 170165
$ws.Range("M89").Value = -164549

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 37000
$ws.Range("I68").Value = 19000
$ws.Range("K68").Value = 19000
$ws.Range("M68").Value = -18251
$ws.Range("H71").Value = 37000
$ws.Range("I71").Value = 19000
$ws.Range("K71").Value = 57000
$ws.Range("M71").Value = -53256
$ws.Range("H86").Value = 10002
$ws.Range("I86").Value = 10002
$ws.Range("K86").Value = 10002
$ws.Range("M86").Value = -8879
$ws.Range("H89").Value = 10002
$ws.Range("I89").Value = 10002
$ws.Range("K89").Value = 50010
$ws.Range("M89").Value = -44394

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 413.5
$ws.Range("J107").Value = 439.25
$ws.Range("L107").Value = 1317.75
$ws.Range("N107").Value = -5157.75
$ws.Range("H131").Value = 4716.6665
$ws.Range("I131").Value = 850
$ws.Range("K131").Value = 2550
$ws.Range("M131").Value = 2490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5474
$ws.Range("I43").Value = 1426.1818
$ws.Range("K43").Value = 1426.1818
$ws.Range("M43").Value = -1275.1818
$ws.Range("H46").Value = 4874.75
$ws.Range("I46").Value = 4874.75
$ws.Range("K46").Value = 4874.75
$ws.Range("M46").Value = -4718.75
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = ""
$ws.Range("H80").Value = 8343.799999999999
$ws.Range("I80").Value = 2222
$ws.Range("J80").Value = 9874.25
$ws.Range("K80").Value = 2222
$ws.Range("L80").Value = 9874.25
$ws.Range("M80").Value = -1224
$ws.Range("N80").Value = -11870.25
$ws.Range("H83").Value = 8343.799999999999
$ws.Range("I83").Value = 2222
$ws.Range("J83").Value = 9874.25
$ws.Range("K83").Value = 11110
$ws.Range("L83").Value = 49371.25
$ws.Range("M83").Value = -6118
$ws.Range("N83").Value = -59355.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = ""
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("H128").Value = 63223.777
$ws.Range("J128").Value = 63223.777
$ws.Range("L128").Value = 63223.777
$ws.Range("N128").Value = -73183.777
$ws.Range("H132").Value = 2412.7856
$ws.Range("I132").Value = 2423.25
$ws.Range("J132").Value = 2350
$ws.Range("K132").Value = 7269.75
$ws.Range("L132").Value = 7050
$ws.Range("M132").Value = -4739.75
$ws.Range("N132").Value = -12110
$ws.Range("H136").Value = 3040.2856
$ws.Range("I136").Value = 3130.3333
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 9390.999899999999
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -6840.999899999999
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6945835
$ws.Range("I126").Value = 6945835
$ws.Range("K126").Value = 20837505
$ws.Range("M126").Value = -20835035
$ws.Range("H132").Value = 142861340
$ws.Range("I132").Value = 4899.3335
$ws.Range("K132").Value = 14698.0005
$ws.Range("M132").Value = -12168.0005
$ws.Range("H136").Value = 8957.200000000001
$ws.Range("I136").Value = 10418.786
$ws.Range("J136").Value = 5546.8335
$ws.Range("K136").Value = 31256.358
$ws.Range("L136").Value = 16640.5005
$ws.Range("M136").Value = -28706.358
$ws.Range("N136").Value = -21740.5005
